# Update cryptos list values (Price and Volume(1h) columns) to reflect latest scrape.
# Cells that would otherwise be auto-converted to numbers by Excel get an explicit
# text NumberFormat first, preserving their original string/text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.559.08'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '2.604.02'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '538.11'
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.70'
$ws.Range("E6").Value = '  +1.72%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.567'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.335'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '3.058.74'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = '59.487.28'
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.75'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").Value = '2.600.51'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.79'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.36'
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.10'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.37'
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.63'
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.24'
$ws.Range("E27").Value = '  +3.09%  '
$ws.Range("D28").Value = '0.0₃0746'
$ws.Range("E28").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.67'
$ws.Range("E30").Value = '  +5.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.82'
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.84'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.98'
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.98'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  -0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.47'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.836'
$ws.Range("E37").Value = '  +2.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.828'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '273.36'
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.600'
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.76'
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0525'
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.67'
$ws.Range("E46").Value = '  +4.19%  '
$ws.Range("D47").Value = '1.943.90'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.51'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.09'
$ws.Range("E50").Value = '  -2.02%  '
$ws.Range("E51").Value = '  +0.49%  '
